$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2808
$ws.Range("J51").Value = 3750.6667
$ws.Range("L51").Value = 3750.6667
$ws.Range("N51").Value = -4718.6667
$ws.Range("H88").Value = 3310.75
$ws.Range("J88").Value = 2498.1428
$ws.Range("L88").Value = 2498.1428
$ws.Range("N88").Value = -3310.1428
$ws.Range("H91").Value = 3310.75
$ws.Range("J91").Value = 2498.1428
$ws.Range("L91").Value = 2498.1428
$ws.Range("N91").Value = -5306.1428
$ws.Range("H106").Value = 2914.1667
$ws.Range("I106").Value = 3475
$ws.Range("J106").Value = 1231.6666
$ws.Range("K106").Value = 3475
$ws.Range("L106").Value = 1231.6666
$ws.Range("M106").Value = -2844
$ws.Range("N106").Value = -2493.6666
$ws.Range("H112").Value = 2174.7144
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 2209.2646
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 6627.793799999999
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -8843.793799999999
$ws.Range("H132").Value = 825.3261
$ws.Range("I132").Value = 764.46344
$ws.Range("K132").Value = 2293.39032
$ws.Range("M132").Value = 236.60968
$ws.Range("H138").Value = 2282.8286
$ws.Range("I138").Value = 2328.3333
$ws.Range("J138").Value = 2234.647
$ws.Range("K138").Value = 6984.999899999999
$ws.Range("L138").Value = 6703.941
$ws.Range("M138").Value = -1844.999899999999
$ws.Range("N138").Value = -16983.941
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5772.6
$ws.Range("I32").Value = 3799.8813
$ws.Range("K32").Value = 3799.8813
$ws.Range("M32").Value = -3512.8813
$ws.Range("H122").Value = 3257
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 3583.1667
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 10749.5001
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -15649.5001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2498.2354
$ws.Range("I20").Value = 2368.7778
$ws.Range("J20").Value = 2643.875
$ws.Range("K20").Value = 2368.7778
$ws.Range("L20").Value = 2643.875
$ws.Range("M20").Value = -2121.7778
$ws.Range("N20").Value = -3137.875
$ws.Range("H86").Value = 146747.86
$ws.Range("I86").Value = 4190
$ws.Range("K86").Value = 4190
$ws.Range("M86").Value = -3067
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H89").Value = 146747.86
$ws.Range("I89").Value = 4190
$ws.Range("K89").Value = 20950
$ws.Range("M89").Value = -15334
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480
$ws.Range("H105").Value = 2353.64
$ws.Range("J105").Value = 4500
$ws.Range("L105").Value = 4500
$ws.Range("N105").Value = -7994
$ws.Range("H107").Value = 1156.7333
$ws.Range("I107").Value = 912.5714
$ws.Range("K107").Value = 912.5714
$ws.Range("M107").Value = 1007.4286
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2139.4
$ws.Range("I122").Value = 1027.6666
$ws.Range("K122").Value = 3082.9998
$ws.Range("M122").Value = -632.9998000000001
$ws.Range("H139").Value = 55780
$ws.Range("J139").Value = 55780
$ws.Range("L139").Value = 55780
$ws.Range("N139").Value = -66060
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 111124180
$ws.Range("I17").Value = 1000000000
$ws.Range("K17").Value = 3000000000
$ws.Range("M17").Value = -2999999831
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").ClearContents()
$ws.Range("H96").Value = 5999.75
$ws.Range("J96").Value = 5999.75
$ws.Range("L96").Value = 17999.25
$ws.Range("N96").Value = -22117.25
$ws.Range("H105").Value = 5399.7
$ws.Range("J105").Value = 5399.7
$ws.Range("L105").Value = 16199.1
$ws.Range("N105").Value = -21441.1
$ws.Range("H131").Value = 15612.681
$ws.Range("J131").Value = 16642.432
$ws.Range("L131").Value = 49927.296
$ws.Range("N131").Value = -60007.296
$ws.Range("H139").Value = 6729
$ws.Range("I139").Value = 6729
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 20187
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -15047
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 1784.1052
$ws.Range("I140").Value = 978.6429000000001
$ws.Range("J140").Value = 4039.4
$ws.Range("K140").Value = 2935.9287
$ws.Range("L140").Value = 12118.2
$ws.Range("M140").Value = 2244.0713
$ws.Range("N140").Value = -22478.2
$ws.Range("H141").Value = 2648.6924
$ws.Range("I141").Value = 2648.6924
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7946.0772
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2766.0772
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4703.6
$ws.Range("I70").Value = 4249.5
$ws.Range("K70").Value = 4249.5
$ws.Range("M70").Value = -3979.5
$ws.Range("H73").Value = 4703.6
$ws.Range("I73").Value = 4249.5
$ws.Range("K73").Value = 4249.5
$ws.Range("M73").Value = -3313.5
$ws.Range("H80").Value = 2383.3684
$ws.Range("I80").Value = 2384.9333
$ws.Range("K80").Value = 2384.9333
$ws.Range("M80").Value = -1386.9333
$ws.Range("H83").Value = 2383.3684
$ws.Range("I83").Value = 2384.9333
$ws.Range("K83").Value = 11924.6665
$ws.Range("M83").Value = -6932.666500000001
$ws.Range("H102").Value = 3214
$ws.Range("I102").Value = 3166.3333
$ws.Range("K102").Value = 3166.3333
$ws.Range("M102").Value = -1544.3333
$ws.Range("H132").Value = 2407027
$ws.Range("I132").Value = 3207469.2
$ws.Range("K132").Value = 9622407.600000001
$ws.Range("M132").Value = -9619877.600000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3269.625
$ws.Range("I7").Value = 3406.1667
$ws.Range("J7").Value = 2860
$ws.Range("K7").Value = 3406.1667
$ws.Range("L7").Value = 2860
$ws.Range("M7").Value = -3294.1667
$ws.Range("N7").Value = -3084
$ws.Range("H94").Value = 39999.5
$ws.Range("J94").Value = 39999.5
$ws.Range("L94").Value = 39999.5
$ws.Range("N94").Value = -41351.5
$ws.Range("H126").Value = 3269.625
$ws.Range("I126").Value = 3406.1667
$ws.Range("J126").Value = 2860
$ws.Range("K126").Value = 10218.5001
$ws.Range("L126").Value = 8580
$ws.Range("M126").Value = -7748.500100000001
$ws.Range("N126").Value = -13520
$ws.Range("H132").Value = 3070.6206
$ws.Range("I132").Value = 1912.2222
$ws.Range("J132").Value = 3591.9
$ws.Range("K132").Value = 5736.6666
$ws.Range("L132").Value = 10775.7
$ws.Range("M132").Value = -3206.6666
$ws.Range("N132").Value = -15835.7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5408.75
$ws.Range("J62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 5408.75
$ws.Range("J65").Value = 7500
$ws.Range("L65").Value = 37500
$ws.Range("N65").Value = -43740
$ws.Range("H81").Value = 2992.077
$ws.Range("I81").Value = 1624.625
$ws.Range("K81").Value = 3249.25
$ws.Range("M81").Value = -2188.25
$ws.Range("H84").Value = 2992.077
$ws.Range("I84").Value = 1624.625
$ws.Range("K84").Value = 16246.25
$ws.Range("M84").Value = -10942.25
$ws.Range("H107").Value = 666.2414
$ws.Range("I107").Value = 522
$ws.Range("K107").Value = 1566
$ws.Range("M107").Value = 354
$ws.Range("H122").Value = 58276.785
$ws.Range("I122").Value = 100120.25
$ws.Range("K122").Value = 300360.75
$ws.Range("M122").Value = -297910.75
$ws.Range("H136").Value = 16341227
$ws.Range("J136").Value = 1963.5834
$ws.Range("L136").Value = 5890.7502
$ws.Range("N136").Value = -10990.7502
